# Update the salary header "Luong" -> "Luong Thang" and refresh the
# view state (column width autofit + active selection) to match the
# latest export from the source Excel file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# E1 header text: "Lương" -> "Lương Tháng"
$ws.Range("E1").Value = "Lương Tháng"

# The new, longer header text no longer fits the default column width,
# so size column E to fit its contents (as Excel does automatically).
$ws.Columns.Item(5).AutoFit()

# Move/restore the active selection to G7, like in the saved file.
$ws.Range("G7").Select()
